# Update column G ("K") values on Sheet1 as part of regenerating save_data
# to use K (strikeouts) instead of the old Strike# stat, and recalculated
# std/mean based s_vals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 3
    3  = 5
    4  = 3
    5  = 9
    6  = 6
    7  = 6
    8  = 4
    9  = 4
    10 = 7
    11 = 8
    12 = 8
    13 = 1
    14 = 6
    15 = 3
    16 = 8
    17 = 4
    18 = 4
    19 = 8
    20 = 6
    21 = 2
    22 = 5
    23 = 5
    24 = 5
    25 = 4
    26 = 5
    27 = 3
    28 = 4
    29 = 5
    30 = 5
    31 = 4
    32 = 5
    33 = 2
    34 = 4
    35 = 5
    36 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
